# Commit: Updated symbol list on Sun Jan 15 16:45:51 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    # A leading apostrophe forces Excel to store the entry as literal text
    # (matching the existing inline-string cells) instead of silently
    # auto-converting numeric-/percent-looking strings into real numbers,
    # which would lose formatting such as trailing zeros ("0.3260", "4.620").
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    # Typing a quote-prefixed value tags the cell with a "quote prefix"
    # style; resetting to Normal keeps the cell's style identical to before.
    $range.Style = "Normal"
}

Set-TextValue "D2" "300.86"
Set-TextValue "E2" "-1.03%"
Set-TextValue "E3" "-3.67%"
Set-TextValue "D4" "5.161"
Set-TextValue "E4" "-2.51%"
Set-TextValue "D5" "0.07378"
Set-TextValue "E5" "-1.46%"
Set-TextValue "D6" "2.105"
Set-TextValue "E6" "41.44%"
Set-TextValue "D7" "7.894"
Set-TextValue "E7" "0.89%"
Set-TextValue "D8" "3.762"
Set-TextValue "E8" "-1.01%"
Set-TextValue "D9" "0.9285"
Set-TextValue "D10" "0.1709"
Set-TextValue "E10" "1.07%"
Set-TextValue "D11" "0.07478"
Set-TextValue "E11" "-4.63%"
Set-TextValue "D12" "0.08181"
Set-TextValue "E12" "2.14%"
Set-TextValue "D13" "0.03038"
Set-TextValue "E13" "0.19%"
Set-TextValue "D14" "0.09918"
Set-TextValue "E14" "0.15%"
Set-TextValue "D15" "0.001493"
Set-TextValue "E15" "0.12%"
Set-TextValue "D16" "0.006101"
Set-TextValue "E16" "-1.66%"
Set-TextValue "D17" "3.448"
Set-TextValue "E17" "-0.73%"
Set-TextValue "D18" "2.228"
Set-TextValue "E18" "-0.08%"
Set-TextValue "D19" "0.3260"
Set-TextValue "E19" "-2.02%"
Set-TextValue "D20" "0.1338"
Set-TextValue "E20" "-0.52%"
Set-TextValue "D21" "4.620"
Set-TextValue "E21" "3.21%"
Set-TextValue "D22" "0.04654"
Set-TextValue "E22" "0.78%"
Set-TextValue "D23" "0.1583"
Set-TextValue "E23" "-2.25%"
Set-TextValue "E24" "-0.07%"
Set-TextValue "E25" "0.88%"
Set-TextValue "E26" "-7.18%"
Set-TextValue "E27" "7.47%"
Set-TextValue "D39" "0.01719"
Set-TextValue "E39" "-1.80%"
Set-TextValue "D40" "0.04511"
Set-TextValue "E40" "-0.72%"
Set-TextValue "D41" "0.007104"
Set-TextValue "E41" "-0.96%"
Set-TextValue "D42" "0.1343"
Set-TextValue "E42" "-0.10%"
Set-TextValue "D43" "0.002128"
Set-TextValue "E43" "-3.66%"
Set-TextValue "D44" "0.01048"
Set-TextValue "E44" "-17.07%"
Set-TextValue "D45" "0.00006271"
Set-TextValue "E45" "1.88%"
Set-TextValue "D46" "0.006995"
Set-TextValue "D47" "1.849"
Set-TextValue "E47" "161.14%"
